$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.351.66"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").Value = "3.376.64"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'573.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("D6").Value = "'136.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "3.376.30"
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("E9").Value = "  -0.34%  "

# Row 10
$ws.Range("D10").Value = "'7.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.83%  "

# Row 11
$ws.Range("E11").Value = "  +1.18%  "

# Row 12
$ws.Range("E12").Value = "  -0.50%  "

# Row 13
$ws.Range("D13").Value = "3.953.72"
$ws.Range("E13").Value = "  -0.11%  "

# Row 14
$ws.Range("E14").Value = "  +2.47%  "

# Row 15
$ws.Range("E15").Value = "  +1.64%  "

# Row 16
$ws.Range("D16").Value = "'25.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.65%  "

# Row 17
$ws.Range("D17").Value = "3.377.97"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18
$ws.Range("D18").Value = "61.429.00"
$ws.Range("E18").Value = "  -0.06%  "

# Row 19
$ws.Range("E19").Value = "  +0.18%  "

# Row 20
$ws.Range("D20").Value = "'5.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.81%  "

# Row 21
$ws.Range("D21").Value = "'9.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.55%  "

# Row 22
$ws.Range("D22").Value = "'376.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.46%  "

# Row 23
$ws.Range("E23").Value = "  -3.25%  "

# Row 24
$ws.Range("D24").Value = "3.515.28"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("E25").Value = "  -0.19%  "

# Row 26
$ws.Range("D26").Value = "'0.0000127"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.41%  "

# Row 27
$ws.Range("D27").Value = "'71.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "

# Row 28
$ws.Range("D28").Value = "'1.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.23%  "

# Row 29
$ws.Range("D29").Value = "'7.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.71%  "

# Row 30
$ws.Range("E30").Value = "  +0.41%  "

# Row 31
$ws.Range("D31").Value = "'8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.47%  "

# Row 32
$ws.Range("E32").Value = "  +2.43%  "

# Row 33
$ws.Range("E33").Value = "  +1.03%  "

# Row 34
$ws.Range("E34").Value = "  +0.04%  "

# Row 35
$ws.Range("D35").Value = "'23.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "

# Row 36
$ws.Range("D36").Value = "'5.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.98%  "

# Row 37
$ws.Range("D37").Value = "'6.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.22%  "

# Row 38
$ws.Range("E38").Value = "  -1.36%  "

# Row 39
$ws.Range("D39").Value = "'165.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.96%  "

# Row 40
$ws.Range("E40").Value = "  -3.18%  "

# Row 41
$ws.Range("E41").Value = "  -0.10%  "

# Row 42
$ws.Range("E42").Value = "  +5.85%  "

# Row 43
$ws.Range("D43").Value = "'0.774"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "

# Row 44
$ws.Range("E44").Value = "  +0.17%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.34%  "

# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'41.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.32%  "

# Row 47
$ws.Range("D47").Value = "'24.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.20%  "

# Row 48
$ws.Range("D48").Value = "'6.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.57%  "

# Row 49
$ws.Range("D49").Value = "'22.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.20%  "

# Row 50
$ws.Range("D50").Value = "2.347.30"
$ws.Range("E50").Value = "  +3.88%  "

# Row 51
$ws.Range("E51").Value = "  +0.61%  "
